# Fruta / hortaliza, semanal
# Inserts two new weekly price-report rows (Damasco, Vega Modelo de Temuco)
# right before the current row 25, pushing the existing rows 25-78 down to
# rows 27-80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 25 (existing rows 25:78 shift down to 27:80)
$ws.Rows("25:26").Insert()

# ---- New row 25 ----
$ws.Cells.Item(25, 1).Value = 10
$ws.Cells.Item(25, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value = "La Araucanía"
$ws.Cells.Item(25, 4).Value = 44914
$ws.Cells.Item(25, 5).Value = 9
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100103
$ws.Cells.Item(25, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(25, 9).Value = 100103003
$ws.Cells.Item(25, 10).Value = "Damasco"
$ws.Cells.Item(25, 11).Value = "Castle Brite"
$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 20
$ws.Cells.Item(25, 14).Value = 19000
$ws.Cells.Item(25, 15).Value = 19000
$ws.Cells.Item(25, 16).Value = 19000
$ws.Cells.Item(25, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(25, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(25, 19).Value = 1900
$ws.Cells.Item(25, 20).Value = 10

# ---- New row 26 ----
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 44914
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100103
$ws.Cells.Item(26, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(26, 9).Value = 100103003
$ws.Cells.Item(26, 10).Value = "Damasco"
$ws.Cells.Item(26, 11).Value = "Castle Brite"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 45
$ws.Cells.Item(26, 14).Value = 16000
$ws.Cells.Item(26, 15).Value = 16000
$ws.Cells.Item(26, 16).Value = 16000
$ws.Cells.Item(26, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(26, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(26, 19).Value = 1600
$ws.Cells.Item(26, 20).Value = 10
